$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 0.79821795875583834
$ws.Range("AJ1").Value = 0.92527283090863821
$ws.Range("P3").Value = 0.94802636237312088
$ws.Range("AS3").Value = 0.64938771677260665
$ws.Range("AZ3").Value = 0.86861671775433957
$ws.Range("B4").Value = 0.66879315422556562
$ws.Range("AT4").Value = 0.75248437145777236
$ws.Range("G5").Value = 0.91244217088861534
$ws.Range("BF5").Value = 0.81814286304133066
$ws.Range("G6").Value = 0.84587132537764509
$ws.Range("BB6").Value = 0.69658890160857534
$ws.Range("A7").Value = 0.77910382604561423
$ws.Range("AI7").Value = 0.79904117880379033
$ws.Range("V8").Value = 0.65698878398041172
$ws.Range("T9").Value = 0.66780914533133751
$ws.Range("V9").Value = 0.94166600014554247
$ws.Range("Y9").Value = 0.72194635117977901
$ws.Range("AW9").Value = 0.58481442429456787
$ws.Range("H10").Value = 0.71754416546673616
$ws.Range("BE10").Value = 0.88626283001362194
$ws.Range("AA13").Value = 0.6785520156863839
$ws.Range("BC13").Value = 0.76553982743486992
$ws.Range("F14").Value = 0.9602431727502212
$ws.Range("AN14").Value = 0.93762744884652927
$ws.Range("BB14").Value = 0.78716386534539007
$ws.Range("BI14").Value = 0.92677922592639783
$ws.Range("Z15").Value = 0.97266863341173648
$ws.Range("AW16").Value = 0.92804730983762762
$ws.Range("Q18").Value = 0.80156808810832003
$ws.Range("AY18").Value = 0.92738396771081666
$ws.Range("D19").Value = 0.9446182105821137
$ws.Range("AF19").Value = 0.86103674968705191
$ws.Range("BE19").Value = 0.98346371869301696
$ws.Range("BL19").Value = 0.91461022716690121
$ws.Range("R20").Value = 0.66961584785577866
$ws.Range("S20").Value = 0.74778916937774054
$ws.Range("AK20").Value = 0.63634255581883781
$ws.Range("BJ20").Value = 0.79736567054333485
$ws.Range("E21").Value = 0.80298655549117748
$ws.Range("W21").Value = 0.62065155854744036
$ws.Range("BF21").Value = 0.84533568372017498
$ws.Range("BM21").Value = 0.91538445771386889
$ws.Range("U22").Value = 0.90352773981993195
$ws.Range("AH23").Value = 0.8966418369082374
$ws.Range("AN23").Value = 0.94710999023757281
$ws.Range("BC24").Value = 0.867308775790814
$ws.Range("B26").Value = 0.82186823835622469
$ws.Range("X26").Value = 0.99459287298554755
$ws.Range("Y26").Value = 0.87587032138079857
$ws.Range("AB26").Value = 0.93719901638178671
$ws.Range("AC26").Value = 0.81571478461373026
$ws.Range("C27").Value = 0.81646313486198419
$ws.Range("I27").Value = 0.84566039907982815
$ws.Range("J28").Value = 0.86307542078166022
$ws.Range("AC28").Value = 0.88325396398164835
$ws.Range("AD29").Value = 0.94101814828648012
$ws.Range("AK29").Value = 0.91738854488452559
$ws.Range("J30").Value = 0.98806886116885073
$ws.Range("AE30").Value = 0.90340882383133003
$ws.Range("AO30").Value = 0.92838675553063821
$ws.Range("AC31").Value = 0.99561810413072238
$ws.Range("BH31").Value = 0.67155002507269912
$ws.Range("BM31").Value = 0.71142235627810713
$ws.Range("BN31").Value = 0.87101848419832706
$ws.Range("G32").Value = 0.82630594711146754
$ws.Range("AR32").Value = 0.9860418628181522
$ws.Range("BA32").Value = 0.62954267621861815
$ws.Range("Q33").Value = 0.82706457591667215
$ws.Range("AF33").Value = 0.60308239041128553
$ws.Range("L35").Value = 0.58206333809575539
$ws.Range("Y35").Value = 0.98684510570260864
$ws.Range("AO35").Value = 0.91719757772662613
$ws.Range("AI37").Value = 0.90896725685455193
$ws.Range("AL37").Value = 0.80159542233087033
$ws.Range("BC37").Value = 0.61849411389349607
$ws.Range("K38").Value = 0.70548555134855295
$ws.Range("AJ38").Value = 0.90147252125373223
$ws.Range("AM38").Value = 0.77627868946262435
$ws.Range("AV38").Value = 0.98253827823759665
$ws.Range("X39").Value = 0.93048440778286945
$ws.Range("Z39").Value = 0.97690861626913816
$ws.Range("Z40").Value = 0.89859352337577669
$ws.Range("O41").Value = 0.65685339796006437
$ws.Range("AS41").Value = 0.97345555637201042
$ws.Range("BL41").Value = 0.71415382641435277
$ws.Range("AJ42").Value = 0.73003953403300548
$ws.Range("BH42").Value = 0.82781347503698566
$ws.Range("W43").Value = 0.62140374393625364
$ws.Range("AG43").Value = 0.95325302720306626
$ws.Range("BP43").Value = 0.59231085191740751
$ws.Range("K44").Value = 0.94497689187669542
$ws.Range("AY44").Value = 0.71633888074762586
$ws.Range("BJ44").Value = 0.87107164787796054
$ws.Range("R45").Value = 0.94341809453141101
$ws.Range("AJ45").Value = 0.84148261273760738
$ws.Range("AX46").Value = 0.85894631338683647
$ws.Range("AZ46").Value = 0.95858329347185611
$ws.Range("BP46").Value = 0.97562921338045161
$ws.Range("L47").Value = 0.91908496928317995
$ws.Range("Y47").Value = 0.68193710505857319
$ws.Range("AW47").Value = 0.83305161539437989
$ws.Range("BM47").Value = 0.97298168344835956
$ws.Range("S48").Value = 0.60341256643472951
$ws.Range("AW48").Value = 0.72946832138116502
$ws.Range("AX49").Value = 0.86951109361394896
$ws.Range("AZ50").Value = 0.68999804517286012
$ws.Range("AG51").Value = 0.94212685282173136
$ws.Range("BJ52").Value = 0.98783252946903644
$ws.Range("J53").Value = 0.79356369734218446
$ws.Range("AU53").Value = 0.99689106829249541
$ws.Range("AV55").Value = 0.61690955152643856
$ws.Range("Q56").Value = 0.66712571696183187
$ws.Range("AB56").Value = 0.88016019699690606
$ws.Range("AH57").Value = 0.90746782814697668
$ws.Range("AV57").Value = 0.71875368200701262
$ws.Range("Q59").Value = 0.73252604703687396
$ws.Range("AF60").Value = 0.89875189737676786
$ws.Range("BG60").Value = 0.78542385379758506
$ws.Range("AS61").Value = 0.85307023021045514
$ws.Range("BK62").Value = 0.94428439653765839
$ws.Range("BB63").Value = 0.85621067773431192
$ws.Range("BG63").Value = 0.73829008178588085
$ws.Range("BM63").Value = 0.78020284462989276
$ws.Range("H64").Value = 0.97158818074988673
$ws.Range("AH64").Value = 0.58671147126266576
$ws.Range("AZ64").Value = 0.98870429235093982
$ws.Range("B65").Value = 0.89756391939216718
$ws.Range("AT65").Value = 0.92758884428291455
$ws.Range("P66").Value = 0.97622662166870722
$ws.Range("AK66").Value = 0.91371290123570448
$ws.Range("H67").Value = 0.89894470229685386
$ws.Range("BG67").Value = 0.85961276446261015
$ws.Range("BH67").Value = 0.76520282542070972
$ws.Range("BL67").Value = 0.8841778237958311
$ws.Range("Z68").Value = 0.98209296783129607
$ws.Range("AK68").Value = 0.91671813760236975
